$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 3967.4583
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 3967.4583
$ws.Range("K17").Value = 0
$ws.Range("M17").Value = 11902.3749
$ws.Range("N17").Value = -12238.3749
$ws.Range("L17").ClearContents()
$ws.Range("H19").Value = 664
$ws.Range("I19").Value = 622.2857
$ws.Range("J19").Value = 737
$ws.Range("K19").Value = 622.2857
$ws.Range("L19").Value = 737
$ws.Range("M19").Value = -447.2857
$ws.Range("N19").Value = -1087
$ws.Range("H32").Value = 6519.5713
$ws.Range("J32").Value = 6439.5
$ws.Range("L32").Value = 6439.5
$ws.Range("N32").Value = -7091.5
$ws.Range("H69").Value = 37637.547
$ws.Range("I69").Value = 220013
$ws.Range("J69").Value = 19400
$ws.Range("K69").Value = 660039
$ws.Range("L69").Value = 58200
$ws.Range("M69").Value = -659165
$ws.Range("N69").Value = -59948
$ws.Range("H72").Value = 37637.547
$ws.Range("I72").Value = 220013
$ws.Range("J72").Value = 19400
$ws.Range("K72").Value = 1980117
$ws.Range("L72").Value = 174600
$ws.Range("M72").Value = -1975749
$ws.Range("N72").Value = -183336
$ws.Range("H125").Value = 300000830
$ws.Range("J125").Value = 166667840
$ws.Range("L125").Value = 1500010560
$ws.Range("N125").Value = -1500015480

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9596.25
$ws.Range("I32").Value = 9596.25
$ws.Range("K32").Value = 9596.25
$ws.Range("M32").Value = -9309.25
$ws.Range("H45").Value = 8665.637000000001
$ws.Range("I45").Value = 8552.200000000001
$ws.Range("K45").Value = 8552.200000000001
$ws.Range("M45").Value = -8175.200000000001
$ws.Range("H61").Value = 2908.5278
$ws.Range("I61").Value = 1256.1111
$ws.Range("K61").Value = 1256.1111
$ws.Range("M61").Value = -1044.1111
$ws.Range("H88").Value = 666.3333
$ws.Range("I88").Value = 699.5
$ws.Range("J88").Value = 600
$ws.Range("K88").Value = 699.5
$ws.Range("L88").Value = 600
$ws.Range("M88").Value = -293.5
$ws.Range("N88").Value = -1412
$ws.Range("H91").Value = 666.3333
$ws.Range("I91").Value = 699.5
$ws.Range("J91").Value = 600
$ws.Range("K91").Value = 699.5
$ws.Range("L91").Value = 600
$ws.Range("M91").Value = 704.5
$ws.Range("N91").Value = -3408
$ws.Range("H132").Value = 2374.75
$ws.Range("J132").Value = 3000
$ws.Range("L132").Value = 9000
$ws.Range("N132").Value = -14060
$ws.Range("H136").Value = 2908.5278
$ws.Range("I136").Value = 1256.1111
$ws.Range("K136").Value = 3768.3333
$ws.Range("M136").Value = -1218.3333

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 507.64
$ws.Range("J80").Value = 197.61539
$ws.Range("L80").Value = 197.61539
$ws.Range("N80").Value = -2193.61539
$ws.Range("H83").Value = 507.64
$ws.Range("J83").Value = 197.61539
$ws.Range("L83").Value = 988.0769499999999
$ws.Range("N83").Value = -10972.07695
$ws.Range("H86").Value = 1766.5
$ws.Range("I86").Value = 1524.75
$ws.Range("J86").Value = 2250
$ws.Range("K86").Value = 1524.75
$ws.Range("L86").Value = 2250
$ws.Range("M86").Value = -401.75
$ws.Range("N86").Value = -4496
$ws.Range("H89").Value = 1766.5
$ws.Range("I89").Value = 1524.75
$ws.Range("J89").Value = 2250
$ws.Range("K89").Value = 7623.75
$ws.Range("L89").Value = 11250
$ws.Range("M89").Value = -2007.75
$ws.Range("N89").Value = -22482
$ws.Range("H105").Value = 2061.4
$ws.Range("I105").Value = 2203.75
$ws.Range("K105").Value = 2203.75
$ws.Range("M105").Value = -456.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 1887049.1
$ws.Range("I6").Value = 2201232.2
$ws.Range("J6").Value = 1950
$ws.Range("K6").Value = 2201232.2
$ws.Range("L6").Value = 1950
$ws.Range("M6").Value = -2201119.2
$ws.Range("N6").Value = -2176
$ws.Range("H31").Value = 2044.2778
$ws.Range("I31").Value = 1685.3334
$ws.Range("J31").Value = 2223.75
$ws.Range("K31").Value = 1685.3334
$ws.Range("L31").Value = 2223.75
$ws.Range("M31").Value = -1390.3334
$ws.Range("N31").Value = -2813.75
$ws.Range("H34").Value = 2044.2778
$ws.Range("I34").Value = 1685.3334
$ws.Range("J34").Value = 2223.75
$ws.Range("K34").Value = 1685.3334
$ws.Range("L34").Value = 2223.75
$ws.Range("M34").Value = -1483.3334
$ws.Range("N34").Value = -2627.75
$ws.Range("H58").Value = 1399.8148
$ws.Range("I58").Value = 1319.0385
$ws.Range("K58").Value = 1319.0385
$ws.Range("M58").Value = -1116.0385
$ws.Range("H86").Value = 9570.299999999999
$ws.Range("I86").Value = 9463.125
$ws.Range("K86").Value = 9463.125
$ws.Range("M86").Value = -8340.125
$ws.Range("H89").Value = 9570.299999999999
$ws.Range("I89").Value = 9463.125
$ws.Range("K89").Value = 47315.625
$ws.Range("M89").Value = -41699.625
$ws.Range("H132").Value = 1899
$ws.Range("I132").Value = 1768
$ws.Range("K132").Value = 5304
$ws.Range("M132").Value = -2774
$ws.Range("H134").Value = 2556
$ws.Range("I134").Value = 2828.1
$ws.Range("J134").Value = 1649
$ws.Range("K134").Value = 8484.299999999999
$ws.Range("L134").Value = 4947
$ws.Range("M134").Value = -5949.299999999999
$ws.Range("N134").Value = -10017
$ws.Range("H136").Value = 1399.8148
$ws.Range("I136").Value = 1319.0385
$ws.Range("K136").Value = 3957.1155
$ws.Range("M136").Value = -1407.1155

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2442.375
$ws.Range("I5").Value = 1457.5
$ws.Range("K5").Value = 4372.5
$ws.Range("M5").Value = -4260.5
$ws.Range("H56").Value = 11212.846
$ws.Range("I56").Value = 11212.846
$ws.Range("K56").Value = 11212.846
$ws.Range("M56").Value = -10682.846
$ws.Range("H131").Value = 558356.1
$ws.Range("J131").Value = 591083.25
$ws.Range("L131").Value = 1773249.75
$ws.Range("N131").Value = -1783329.75
$ws.Range("H132").Value = 2528.8333
$ws.Range("I132").Value = 585
$ws.Range("J132").Value = 2917.6
$ws.Range("K132").Value = 5265
$ws.Range("L132").Value = 26258.4
$ws.Range("M132").Value = -2735
$ws.Range("N132").Value = -31318.4
$ws.Range("H135").Value = 2442.375
$ws.Range("I135").Value = 1457.5
$ws.Range("K135").Value = 13117.5
$ws.Range("M135").Value = -10582.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 3140205.5
$ws.Range("I3").Value = 3140205.5
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 3140205.5
$ws.Range("L3").Value = 0
$ws.Range("N3").Value = -3140089.5
$ws.Range("M3").ClearContents()
$ws.Range("H10").Value = 4420.6
$ws.Range("I10").Value = 368
$ws.Range("J10").Value = 10499.5
$ws.Range("K10").Value = 368
$ws.Range("L10").Value = 10499.5
$ws.Range("M10").Value = -199
$ws.Range("N10").Value = -10837.5
$ws.Range("H29").Value = 12300000
$ws.Range("I29").Value = 12300000
$ws.Range("K29").Value = 12300000
$ws.Range("M29").Value = -12299710

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3443.75
$ws.Range("H82").Value = 1174.75
$ws.Range("I82").Value = 1266
$ws.Range("K82").Value = 1266
$ws.Range("M82").Value = -905
$ws.Range("H85").Value = 1174.75
$ws.Range("I85").Value = 1266
$ws.Range("K85").Value = 1266
$ws.Range("M85").Value = -18
$ws.Range("H98").Value = 20000
$ws.Range("J98").Value = 20000
$ws.Range("L98").Value = 20000
$ws.Range("N98").Value = -25990
$ws.Range("H100").Value = 4005.182
$ws.Range("J100").Value = 5499.3335
$ws.Range("L100").Value = 5499.3335
$ws.Range("N100").Value = -6581.3335

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1998.2
$ws.Range("I107").Value = 1998.25
$ws.Range("J107").Value = 1998
$ws.Range("K107").Value = 5994.75
$ws.Range("L107").Value = 5994
$ws.Range("M107").Value = -4074.75
$ws.Range("N107").Value = -9834
$ws.Range("H126").Value = 4080.25
$ws.Range("I126").Value = 2492.3845
$ws.Range("J126").Value = 7029.143
$ws.Range("K126").Value = 7477.1535
$ws.Range("L126").Value = 21087.429
$ws.Range("M126").Value = -5007.1535
$ws.Range("N126").Value = -26027.429
$ws.Range("H132").Value = 5066.3125
$ws.Range("I132").Value = 6205.727
$ws.Range("K132").Value = 18617.181
$ws.Range("M132").Value = -16087.181
